$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values only
$ws.Range("B3").Value = 0.9981734024742205
$ws.Range("C3").Value = 0.9981728679661583
$ws.Range("D3").Value = 0.9963629770444413

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9981207691050202
$ws.Range("C4").Value = 0.9982289403384877
$ws.Range("D4").Value = 0.9963231653004051

# Row 5: AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9986668374692407
$ws.Range("C5").Value = 0.9985228570479242
$ws.Range("D5").Value = 0.998834272629967
